# changed OCCU variables, added set.seed to bootstrapping
# Updated computed results (bootstrap estimates/CIs) on both sheets
# to reflect the re-run with a fixed random seed.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Full results")
$ws.Range("C2").Value = 0.646323325211804
$ws.Range("D2").Value = 0.353901374928943
$ws.Range("E2").Value = 1.00022470014075
$ws.Range("J2").Value = 0.353821871104706
$ws.Range("K2").Value = 0.306030320695031
$ws.Range("L2").Value = 0.0140135097358763
$ws.Range("M2").Value = 0.044042949480424
$ws.Range("N2").Value = 0.320043830430907
$ws.Range("F3").Value = 0.61628713784792
$ws.Range("G3").Value = 0.306099085751164
$ws.Range("H4").Value = 0.602270479274433
$ws.Range("I4").Value = 0.285664744666187
$ws.Range("O4").Value = 0.39786482058513
$ws.Range("C5").Value = 0.87736090521989
$ws.Range("D5").Value = 0.122717006372813
$ws.Range("E5").Value = 1.0000779115927
$ws.Range("J5").Value = 0.122707446040255
$ws.Range("K5").Value = 0.124661493700764
$ws.Range("L5").Value = 0.0258778380079757
$ws.Range("M5").Value = 0.0575515385189966
$ws.Range("N5").Value = 0.15053933170874
$ws.Range("F6").Value = 0.845684736960415
$ws.Range("G6").Value = 0.124671206276287
$ws.Range("H7").Value = 0.819804882768865
$ws.Range("I7").Value = 0.125102165595073
$ws.Range("O7").Value = 0.180258984559252
$ws.Range("C8").Value = 0.835673869589237
$ws.Range("D8").Value = 0.164430530342826
$ws.Range("E8").Value = 1.00010439993206
$ws.Range("J8").Value = 0.164413365598627
$ws.Range("K8").Value = 0.155149501851177
$ws.Range("L8").Value = 0.0418768879942465
$ws.Range("M8").Value = 0.0497546069970125
$ws.Range("N8").Value = 0.197026389845424
$ws.Range("F9").Value = 0.827795328153142
$ws.Range("G9").Value = 0.15516569944863
$ws.Range("H10").Value = 0.785914068214634
$ws.Range("I10").Value = 0.129262717242758
$ws.Range("O10").Value = 0.21416797259564
$ws.Range("C11").Value = 0.768779677316853
$ws.Range("D11").Value = 0.231367222074937
$ws.Range("E11").Value = 1.00014689939179
$ws.Range("J11").Value = 0.231333239362774
$ws.Range("K11").Value = 0.225387262095531
$ws.Range("L11").Value = 0.0283411036012153
$ws.Range("M11").Value = 0.0328616107307432
$ws.Range("N11").Value = 0.253728365696746
$ws.Range("F12").Value = 0.764258506127577
$ws.Range("G12").Value = 0.22542037134725
$ws.Range("H13").Value = 0.73591323923548
$ws.Range("I13").Value = 0.220052021730966
$ws.Range("O13").Value = 0.264194850093517
$ws.Range("C14").Value = 0.86056319164205
$ws.Range("D14").Value = 0.139525364609732
$ws.Range("E14").Value = 1.00008855625178
$ws.Range("J14").Value = 0.139513009860504
$ws.Range("K14").Value = 0.12889266239111
$ws.Range("L14").Value = 0.0156862351653176
$ws.Range("M14").Value = 0.0268614147365685
$ws.Range("N14").Value = 0.144578897556428
$ws.Range("F15").Value = 0.849387022438783
$ws.Range("G15").Value = 0.128904076642174
$ws.Range("H16").Value = 0.833699398159275
$ws.Range("I16").Value = 0.140850119093156
$ws.Range("O16").Value = 0.166374424597072

$ws = $wb.Worksheets.Item("For plotting")
$ws.Range("B2").Value = 0.258402143126609
$ws.Range("C2").Value = 0.381685517735206
$ws.Range("D2").Value = 0.320043830430907
$ws.Range("B3").Value = 0.350235918909555
$ws.Range("C3").Value = 0.445493722260706
$ws.Range("D3").Value = 0.39786482058513
$ws.Range("B4").Value = 0.302156483748208
$ws.Range("C4").Value = 0.405487258461205
$ws.Range("D4").Value = 0.353821871104706
$ws.Range("B5").Value = 0.0803137283705362
$ws.Range("C5").Value = 0.220764935046943
$ws.Range("D5").Value = 0.15053933170874
$ws.Range("B6").Value = 0.118786280157073
$ws.Range("C6").Value = 0.241731688961431
$ws.Range("D6").Value = 0.180258984559252
$ws.Range("B7").Value = 0.0639261260280642
$ws.Range("C7").Value = 0.181488766052447
$ws.Range("D7").Value = 0.122707446040255
$ws.Range("B8").Value = 0.0850071193212848
$ws.Range("C8").Value = 0.309045660369563
$ws.Range("D8").Value = 0.197026389845424
$ws.Range("B9").Value = 0.110050467661553
$ws.Range("C9").Value = 0.318285477529726
$ws.Range("D9").Value = 0.21416797259564
$ws.Range("B10").Value = 0.0396943543195528
$ws.Range("C10").Value = 0.289132376877701
$ws.Range("D10").Value = 0.164413365598627
$ws.Range("B11").Value = 0.0952628451618303
$ws.Range("C11").Value = 0.412193886231662
$ws.Range("D11").Value = 0.253728365696746
$ws.Range("B12").Value = 0.116719389489182
$ws.Range("C12").Value = 0.411670310697852
$ws.Range("D12").Value = 0.264194850093517
$ws.Range("B13").Value = 0.0628096520449353
$ws.Range("C13").Value = 0.399856826680612
$ws.Range("D13").Value = 0.231333239362774
$ws.Range("B14").Value = 0.0774956805350221
$ws.Range("C14").Value = 0.211662114577834
$ws.Range("D14").Value = 0.144578897556428
$ws.Range("B15").Value = 0.106213225852818
$ws.Range("C15").Value = 0.226535623341326
$ws.Range("D15").Value = 0.166374424597072
$ws.Range("B16").Value = 0.0764858670750393
$ws.Range("C16").Value = 0.202540152645968
$ws.Range("D16").Value = 0.139513009860504
